$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Row 4 = "Electricity" need row.
# Set electricity use from gas boiler (H4, heating), gas boiler for hot
# sanitary water (K4) and gas stove for cooking (P4) to zero.
$ws.Range("H4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("P4").Value = 0
